# Uren Registratie 04-02-2016 & Bug Report fix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fahrettin (row 6) was absent Monday (grandfather in hospital) and present the rest of the week
$ws.Range("C6").Value2 = 0
$ws.Range("D6").Value2 = 2
$ws.Range("E6").Value2 = 2
$ws.Range("F6").Value2 = 2
$ws.Range("G6").Value2 = 2
$ws.Range("H6").Value2 = "Fahrettin: Opa Ziekenhuis"

# Ruben (row 7) was sick on Wednesday
$ws.Range("C7").Value2 = 8
$ws.Range("D7").Value2 = 8
$ws.Range("E7").Value2 = 0
$ws.Range("F7").Value2 = 8
$ws.Range("G7").Value2 = 7
$ws.Range("H7").Value2 = "Ruben: Ziek"

# Highlight the missed/partial days
$ws.Range("C6").Interior.Color = 255
$ws.Range("E7").Interior.Color = 255

# Update total lesuren (manually entered total) for the week
$ws.Range("B8").Value2 = 18

# Recalculate all dependent formulas
$excel.Calculate()

# Restore last active selection
$ws.Range("K21").Select()
